$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4565.2607
$ws.Range("I32").Value = 5428.8335
$ws.Range("J32").Value = 4260.4707
$ws.Range("K32").Value = 5428.8335
$ws.Range("L32").Value = 4260.4707
$ws.Range("M32").Value = -5102.8335
$ws.Range("N32").Value = -4912.4707
$ws.Range("H98").Value = 2000.4375
$ws.Range("I98").Value = 2083.4285
$ws.Range("K98").Value = 2083.4285
$ws.Range("M98").Value = -585.4285
$ws.Range("H100").Value = 1704.3334
$ws.Range("I100").Value = 1275.25
$ws.Range("K100").Value = 1275.25
$ws.Range("M100").Value = -734.25
$ws.Range("H112").Value = 3294.5557
$ws.Range("J112").Value = 3458.1875
$ws.Range("L112").Value = 10374.5625
$ws.Range("N112").Value = -12590.5625
$ws.Range("H122").Value = 2000.4375
$ws.Range("I122").Value = 2083.4285
$ws.Range("K122").Value = 6250.2855
$ws.Range("M122").Value = -3800.2855
$ws.Range("H129").Value = 1703.6666
$ws.Range("I129").Value = 1444.5
$ws.Range("J129").Value = 2222
$ws.Range("K129").Value = 4333.5
$ws.Range("L129").Value = 6666
$ws.Range("M129").Value = 666.5
$ws.Range("N129").Value = -16666
$ws.Range("H133").Value = 92197
$ws.Range("J133").Value = 92197
$ws.Range("L133").Value = 92197
$ws.Range("N133").Value = -102317

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3818.762
$ws.Range("I45").Value = 1849.5834
$ws.Range("K45").Value = 1849.5834
$ws.Range("M45").Value = -1472.5834
$ws.Range("H46").Value = 6990.091
$ws.Range("J46").Value = 7221.222
$ws.Range("L46").Value = 7221.222
$ws.Range("N46").Value = -7859.222
$ws.Range("H61").Value = 5092
$ws.Range("I61").Value = 4035.8718
$ws.Range("K61").Value = 4035.8718
$ws.Range("M61").Value = -3823.8718
$ws.Range("H74").Value = 4232.5415
$ws.Range("I74").Value = 3796
$ws.Range("K74").Value = 3796
$ws.Range("M74").Value = -2922
$ws.Range("H77").Value = 4232.5415
$ws.Range("I77").Value = 3796
$ws.Range("K77").Value = 18980
$ws.Range("M77").Value = -14612
$ws.Range("H92").Value = 44975
$ws.Range("J92").Value = 44975
$ws.Range("L92").Value = 44975
$ws.Range("N92").Value = -49967
$ws.Range("H97").Value = 823.9375
$ws.Range("I97").Value = 771.61536
$ws.Range("K97").Value = 771.61536
$ws.Range("M97").Value = -275.61536
$ws.Range("H111").Value = 50277.5
$ws.Range("J111").Value = 50277.5
$ws.Range("L111").Value = 50277.5
$ws.Range("N111").Value = -58457.5
$ws.Range("H122").Value = 3494.25
$ws.Range("I122").Value = 3494.25
$ws.Range("K122").Value = 10482.75
$ws.Range("M122").Value = -8032.75
$ws.Range("H136").Value = 5092
$ws.Range("I136").Value = 4035.8718
$ws.Range("K136").Value = 12107.6154
$ws.Range("M136").Value = -9557.615399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 100
$ws.Range("J10").Value = 100
$ws.Range("L10").Value = 100
$ws.Range("N10").Value = -380
$ws.Range("H22").Value = 4600
$ws.Range("I22").Value = 3500
$ws.Range("K22").Value = 3500
$ws.Range("M22").Value = -3327
$ws.Range("H76").Value = 33333
$ws.Range("J76").Value = 33333
$ws.Range("L76").Value = 33333
$ws.Range("H79").Value = 33333
$ws.Range("J79").Value = 33333
$ws.Range("L79").Value = 33333
$ws.Range("H80").Value = 754
$ws.Range("J80").Value = 838.6667
$ws.Range("L80").Value = 838.6667
$ws.Range("N80").Value = -2834.6667
$ws.Range("H83").Value = 754
$ws.Range("J83").Value = 838.6667
$ws.Range("L83").Value = 4193.3335
$ws.Range("N83").Value = -14177.3335
$ws.Range("H99").Value = 5409.5884
$ws.Range("I99").Value = 4136.4
$ws.Range("K99").Value = 4136.4
$ws.Range("M99").Value = -2638.4
$ws.Range("H106").Value = 14999
$ws.Range("J106").Value = 14999
$ws.Range("L106").Value = 14999
$ws.Range("N106").Value = -17523
$ws.Range("N76").Value = -33963
$ws.Range("N79").Value = -35517

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 104.125
$ws.Range("I7").Value = 105.666664
$ws.Range("K7").Value = 105.666664
$ws.Range("M7").Value = 7.333336000000003
$ws.Range("H140").Value = 114242.2
$ws.Range("J140").Value = 114242.2
$ws.Range("L140").Value = 114242.2
$ws.Range("N140").Value = -124602.2
$ws.Range("H141").Value = 31836.4
$ws.Range("I141").Value = 26764.666
$ws.Range("J141").Value = 39444
$ws.Range("K141").Value = 26764.666
$ws.Range("L141").Value = 39444
$ws.Range("M141").Value = -21584.666
$ws.Range("N141").Value = -49804

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 591811.1
$ws.Range("I128").Value = 591811.1
$ws.Range("K128").Value = 1775433.3
$ws.Range("M128").Value = -1770453.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 19850
$ws.Range("J22").Value = 19850
$ws.Range("L22").Value = 19850
$ws.Range("N22").Value = -20908
$ws.Range("H93").Value = 39858.6
$ws.Range("J93").Value = 39858.6
$ws.Range("L93").Value = 39858.6
$ws.Range("N93").Value = -43602.6
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990
$ws.Range("H105").Value = 45000
$ws.Range("J105").Value = 45000
$ws.Range("L105").Value = 45000
$ws.Range("N105").Value = -51988

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3749.5
$ws.Range("J7").Value = 4999.5
$ws.Range("L7").Value = 4999.5
$ws.Range("N7").Value = -5223.5
$ws.Range("H40").Value = 3300.5
$ws.Range("I40").Value = 2517
$ws.Range("K40").Value = 2517
$ws.Range("M40").Value = -2381
$ws.Range("H61").Value = 253839.75
$ws.Range("I61").Value = 253839.75
$ws.Range("K61").Value = 253839.75
$ws.Range("M61").Value = -253637.75
$ws.Range("H68").Value = 2074.5
$ws.Range("I68").Value = 2074.5
$ws.Range("K68").Value = 2074.5
$ws.Range("M68").Value = -1325.5
$ws.Range("H71").Value = 2074.5
$ws.Range("I71").Value = 2074.5
$ws.Range("K71").Value = 10372.5
$ws.Range("M71").Value = -6628.5
$ws.Range("H113").Value = 253839.75
$ws.Range("I113").Value = 253839.75
$ws.Range("K113").Value = 253839.75
$ws.Range("M113").Value = -251669.75
$ws.Range("H126").Value = 3749.5
$ws.Range("J126").Value = 4999.5
$ws.Range("L126").Value = 14998.5
$ws.Range("N126").Value = -19938.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 46999.5
$ws.Range("J104").Value = 46999.5
$ws.Range("L104").Value = 46999.5
$ws.Range("N104").Value = -53987.5
$ws.Range("H122").Value = 10599.6
$ws.Range("I122").Value = 10499
$ws.Range("K122").Value = 31497
$ws.Range("M122").Value = -29047
$ws.Range("H126").Value = 3224.75
$ws.Range("I126").Value = 3299.6667
$ws.Range("K126").Value = 9899.000100000001
$ws.Range("M126").Value = -7429.000100000001
$ws.Range("H136").Value = 5502.1113
$ws.Range("I136").Value = 4972
$ws.Range("J136").Value = 7357.5
$ws.Range("K136").Value = 14916
$ws.Range("L136").Value = 22072.5
$ws.Range("M136").Value = -12366
$ws.Range("N136").Value = -27172.5
